$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.186.77"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.585.98"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'212.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "1.809.17"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.619.48"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "'63.92"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "26.187.60"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "0.0₃0724"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'214.23"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "'8.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").Value = "'144.20"
$ws.Range("D25").ClearFormats()
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'6.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").Value = "1.406.57"
$ws.Range("E33").Value = "  +7.92%  "
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "'0.587"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.57%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "'0.820"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'0.961"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -12.99%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.765"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "1.720.71"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").Value = "'60.96"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("D47").Value = "'85.65"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "'0.0971"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.26%  "

Write-Output "Done applying cryptos update."
